$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3412
$ws.Range("F6").Value = 7956
$ws.Range("F9").Value = 24
$ws.Range("F10").Value = 1918
$ws.Range("F12").Value = 189
$ws.Range("F13").Value = 1804
$ws.Range("F14").Value = 77
$ws.Range("F18").Value = 1105
$ws.Range("F19").Value = 8617
$ws.Range("F20").Value = 220
$ws.Range("F21").Value = 1150
$ws.Range("F22").Value = 320
$ws.Range("F24").Value = 1057
$ws.Range("F25").Value = 1046
$ws.Range("F26").Value = 585
$ws.Range("F27").Value = 1211
$ws.Range("F28").Value = 1083
$ws.Range("F29").Value = 613
$ws.Range("F30").Value = 510
$ws.Range("F32").Value = 1013
$ws.Range("F33").Value = 127
$ws.Range("F34").Value = 1070
$ws.Range("F35").Value = 487
$ws.Range("F37").Value = 3635
$ws.Range("F38").Value = 75
$ws.Range("F39").Value = 47
$ws.Range("F42").Value = 143
$ws.Range("F44").Value = 739
$ws.Range("F45").Value = 72
$ws.Range("F46").Value = 122
$ws.Range("F47").Value = 1010

$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 362
$ws.Range("F10").Value = 210
$ws.Range("F11").Value = 38
$ws.Range("F12").Value = 38
$ws.Range("F25").Value = 7045
$ws.Range("F33").Value = 61
$ws.Range("F34").Value = 10
$ws.Range("F37").Value = 4
$ws.Range("F41").Value = 3

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2116
$ws.Range("F8").Value = 2269
$ws.Range("F9").Value = 9095
$ws.Range("F10").Value = 1384

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 3412
$ws.Range("F4").Value = 2116
$ws.Range("F5").Value = 7956
$ws.Range("F7").Value = 2269
$ws.Range("F8").Value = 1384
$ws.Range("F11").Value = 189
$ws.Range("F12").Value = 1804
$ws.Range("F13").Value = 77
$ws.Range("F17").Value = 1105
$ws.Range("F18").Value = 8617
$ws.Range("F19").Value = 220
$ws.Range("F20").Value = 1150
$ws.Range("F22").Value = 1057
$ws.Range("F23").Value = 1046
$ws.Range("F24").Value = 585
$ws.Range("F25").Value = 1211
$ws.Range("F26").Value = 1083
$ws.Range("F27").Value = 613
$ws.Range("F28").Value = 510
$ws.Range("F29").Value = 1013
$ws.Range("F30").Value = 38
$ws.Range("F32").Value = 127
$ws.Range("F33").Value = 1070
$ws.Range("F34").Value = 487
$ws.Range("F37").Value = 3635
$ws.Range("F38").Value = 75
$ws.Range("F40").Value = 143
$ws.Range("F42").Value = 739
$ws.Range("F44").Value = 72
$ws.Range("F45").Value = 122
$ws.Range("F46").Value = 1010
